$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Update the existing single-use text (row 4) to hold the new placeholder
# text used by the newly created "All Zones" container.
$ws.Range("E4").Value = "<zoneName>"

# Add a new single-use text entry (row 5) for the Zone Name text that is
# displayed inside the newly created container.
$ws.Range("B5").Value = "SingleUseId2"
$ws.Range("C5").Value = "Default"
$ws.Range("D5").Value = "Left"
$ws.Range("E5").Value = "Zone name"
$ws.Range("F5").Value = "LTR"
